$d = $word.ActiveDocument
$d.Content.Find.Execute("Inputting 11+9 will produce the following:", $true, $false, $false, $false, $false, $true, 1, $false, "Inputting 11 +9 will produce the following:", 2)
